# Applies the "Inicio dos experimentos. Planilha atualizada" edit:
#  - Rename sheet "Experimentos 01" -> "Experimentos"
#  - Add a new column A (index column) on the "Experimentos" sheet with
#    header "Experimentos" merged A1:A2, values 1 and 2 in A3/A4
#  - Add a second experiment row (row 4) with SVM / "DEFAULT (with
#    parallel)" and parameters 5 / 100 / 20 (no result yet)
#  - Adjust column widths, selection, dimension, scroll position
#  - Scroll the "Base de dados" sheet so row 4 is at the top

$wb = $excel.ActiveWorkbook

# --- Sheet 2: rename and edit contents ---
$ws2 = $wb.Worksheets.Item("Experimentos 01")
$ws2.Name = "Experimentos"

# New column A header text + centered alignment, then merge A1:A2
$ws2.Range("A1").Value = "Experimentos"
$ws2.Range("A1").HorizontalAlignment = -4108   # xlCenter
$ws2.Range("A1").VerticalAlignment = -4108     # xlCenter
$ws2.Range("A1:A2").MergeCells = $true

# Index values under the new column
$ws2.Range("A3").Value = 1
$ws2.Range("A4").Value = 2

# New experiment row 4
$ws2.Range("B4").Value = "SVM"
$ws2.Range("C4").Value = "DEFAULT (with parallel)"
$ws2.Range("D4").Value = 5
$ws2.Range("E4").Value = 100
$ws2.Range("F4").Value = 20

# Column widths (values chosen so the persisted width matches the target)
$ws2.Columns.Item(1).ColumnWidth = 12.592447916666666   # -> 13.42578125
$ws2.Columns.Item(3).ColumnWidth = 21.307291666666668   # -> 22.140625
$ws2.Columns.Item(5).ColumnWidth = 20.307291666666668   # -> 21.140625
$ws2.Columns.Item(8).ColumnWidth = 11.166666666666666   # -> 12

# Selection: activeCell D4, selection D4:F4
$ws2.Activate()
$ws2.Range("D4:F4").Select()

# --- Sheet 1: "Base de dados" -> scroll the view so row 4 is on top ---
$ws1 = $wb.Worksheets.Item("Base de dados")
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

$ws2.Activate()

$wb.Save()
